$d = $word.ActiveDocument

# 1) "User Login and Logout" -> "The user and admin can login and logout of the website."
$d.Content.Find.Execute("User Login and Logout", $true, $false, $false, $false, $false, $true, 1, $false, `
    "The user and admin can login and logout of the website.", 2) | Out-Null

# 2) "Dashboard of concerns" -> "The user can create concern/s."
$d.Content.Find.Execute("Dashboard of concerns", $true, $false, $false, $false, $false, $true, 1, $false, `
    "The user can create concern/s.", 2) | Out-Null

# 3) "Create concern" -> "The user can view and update their own concern/s."
$d.Content.Find.Execute("Create concern", $true, $false, $false, $false, $false, $true, 1, $false, `
    "The user can view and update their own concern/s.", 2) | Out-Null

# 4) "Review concern" -> "The user can delete their own concern/s."
$d.Content.Find.Execute("Review concern", $true, $false, $false, $false, $false, $true, 1, $false, `
    "The user can delete their own concern/s.", 2) | Out-Null

# 5) "Edit concern" -> "The user can edit their own concern/s."
$d.Content.Find.Execute("Edit concern", $true, $false, $false, $false, $false, $true, 1, $false, `
    "The user can edit their own concern/s.", 2) | Out-Null

# 6) "Delete concern" -> "In the admin dashboard, the admin can monitor and reply to the concerns of the users."
$d.Content.Find.Execute("Delete concern", $true, $false, $false, $false, $false, $true, 1, $false, `
    "In the admin dashboard, the admin can monitor and reply to the concerns of the users.", 2) | Out-Null

# 7) "Show concern log of the student" -> "The admin can create announcements that will be shown in the landing page."
$d.Content.Find.Execute("Show concern log of the student", $true, $false, $false, $false, $false, $true, 1, $false, `
    "The admin can create announcements that will be shown in the landing page.", 2) | Out-Null

# 8) "Rate the service of the system once concern is solved" -> "Once the user has logged in, the user has the option to view his/her concern log."
$d.Content.Find.Execute("Rate the service of the system once concern is solved", $true, $false, $false, $false, $false, $true, 1, $false, `
    "Once the user has logged in, the user has the option to view his/her concern log.", 2) | Out-Null

# 9) Append a brand-new bullet paragraph after it:
#    "The user can rate the service of the system once concern is solved."
$i = 0
$p = $d.Paragraphs.Item(1)
$target = $null
while ($p -ne $null -and $i -lt 200) {
    $i++
    if ($p.Range.Text -like "Once the user has logged in*") {
        $target = $p
        break
    }
    $p = $p.Next()
}

if ($target -ne $null) {
    $r = $target.Range
    $r.Collapse(0)
    $r.InsertAfter("`r")
    $newPara = $target.Next()
    $newRange = $newPara.Range
    $newRange.InsertAfter("The user can rate the service of the system once concern is solved.")
}
